# Regenerate save_data to use K (strikeouts) instead of Strike# placeholder values,
# recomputed std/mean, and write new simulated s_vals into column G (the "K" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K"), rows 2 through 73.
$newG = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 0
    26 = 0
    27 = 2
    28 = 0
    29 = 2
    30 = 0
    31 = 1
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 2
    38 = 2
    39 = 1
    40 = 1
    41 = 0
    42 = 3
    43 = 0
    44 = 2
    45 = 0
    46 = 0
    47 = 1
    48 = 2
    49 = 0
    50 = 2
    51 = 3
    52 = 1
    53 = 1
    54 = 4
    55 = 0
    56 = 0
    57 = 2
    58 = 1
    59 = 0
    60 = 3
    61 = 2
    62 = 0
    63 = 2
    64 = 0
    65 = 0
    66 = 2
    67 = 1
    68 = 1
    69 = 1
    70 = 1
    71 = 1
    72 = 0
    73 = 1
}

foreach ($row in $newG.Keys) {
    $ws.Cells.Item($row, 7).Value = $newG[$row]
}
